$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = '27.123.24'
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = '  +0.65%  '
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = '1.678.51'
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = '  +0.31%  '
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = '  -0.06%  '
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = '''215.05'
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = '  +0.19%  '
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = '''0.518'
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = '  +0.23%  '
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = '  -0.10%  '
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = '  +2.00%  '
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = '''21.51'
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = '  +5.62%  '
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = '''0.0623'
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = '  +0.61%  '
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = '''0.0888'
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = '  +0.16%  '
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = '1.914.58'
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = '  +0.25%  '
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = '1.686.27'
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = '  -0.20%  '
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = '  +1.43%  '
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = '''0.536'
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = '  +1.93%  '
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.Value = '''66.24'
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = '  +0.84%  '
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = '27.115.60'
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = '  +0.54%  '
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = '''238.65'
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = '  +1.09%  '
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = '  +0.42%  '
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = '0.0₃0742'
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = '  +1.33%  '
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = '  -0.05%  '
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = '''4.52'
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = '  +2.04%  '
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = '''9.46'
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = '  +2.96%  '
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = '  -3.86%  '
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.Value = '''147.65'
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = '  +1.55%  '
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = '''7.25'
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = '  +0.23%  '
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = '  +2.04%  '
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = '  +0.49%  '
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = '  +0.05%  '
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = '  +0.11%  '
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = '  +0.13%  '
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.Value = '1.567.52'
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.Value = '''3.37'
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = '  +1.55%  '
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = '''3.21'
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = '  +2.91%  '
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.Value = '''1.69'
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = '  +0.41%  '
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = '''0.600'
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = '  +2.62%  '
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = '  -1.23%  '
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = '''0.933'
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = '  +4.26%  '
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = '''0.0174'
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = '  +1.29%  '
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = '  +2.37%  '
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = '''68.85'
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = '  +2.97%  '
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = '  -0.08%  '
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = '''5.58'
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = '  -5.28%  '
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = '''2.26'
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = '  -2.27%  '
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = '1.822.99'
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = '''0.781'
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = '  +0.84%  '
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = '''90.60'
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = '  +0.18%  '
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = '  +2.98%  '
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = '  +2.13%  '
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = '''8.13'
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = '  +6.13%  '
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = '  +1.76%  '
$c.Style = "Normal"
